# Fix bug in mk_frq_cnts to include entire text node.
#
# The first paragraph of the document is the synthetic "<page>071v</page>"
# marker line. Its first run (rendered in light-grey Courier New) used to
# read "ière<page>" but got truncated down to just "<page>" by a text-node
# bug; this restores the missing leading "ière" as its own run, re-using
# the exact run formatting (Courier New / color a9a9a9 / 9pt / szCs / rtl)
# that the existing "<page>" run already carries.

$d = $word.ActiveDocument

$searchText = "<page>"
$insertText = "ière"

# Scope the search to the first paragraph (where the "<page>" marker lives)
# and locate it from the very start of the document/story.
$scope = $d.Paragraphs(1).Range
$find = $d.Range($scope.Start, $scope.Start)
$found = $find.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text '$searchText'"
}
$originalStart = $find.Start
$originalEnd = $find.End

# Replace "<page>" with "ière<page>" in a single Find/Replace call so the
# freshly-typed "ière" text inherits the exact run-level formatting (font,
# size, color, complex-script size, rtl) of the run it was typed into,
# rather than picking up default/plain formatting.
$find2 = $d.Range($originalStart, $originalEnd)
$find2.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, ($insertText + $searchText), 1) | Out-Null

# At this point "ière" and "<page>" live inside one merged run. Re-assign
# the FormattedText of just the trailing "<page>" portion back onto itself
# so the engine splits the merged run back into two distinct <w:r>
# elements -- one holding "ière", one holding "<page>" -- each keeping the
# full original run-properties (rFonts/color/sz/szCs/rtl) intact, matching
# the structure of a newly-inserted preceding run.
$insertLen = $insertText.Length
$searchLen = $searchText.Length
$splitStart = $originalStart + $insertLen
$splitEnd = $splitStart + $searchLen

$sub = $d.Range($splitStart, $splitEnd)
$sub.FormattedText = $sub.FormattedText

Write-Output "Paragraph1 now reads: [$($d.Paragraphs(1).Range.Text)]"
